# Update the per-nurse summary figures (regular/overtime hours per shift)
# to reflect the results of the updated optimization run using the HiGHS
# solver. Only the data cells in columns B:G for rows 4-33 change; the
# nurse_id column (A) and header rows stay the same.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(7,6,7,0,1,0),
    @(6,7,7,1,0,0),
    @(6,7,7,1,0,1),
    @(6,7,7,1,0,1),
    @(6,7,7,1,0,1),
    @(6,7,7,1,1,0),
    @(7,6,7,0,1,1),
    @(7,7,6,1,1,1),
    @(7,6,7,0,1,0),
    @(7,6,7,0,1,0),
    @(8,5,7,0,2,0),
    @(7,6,7,0,1,0),
    @(7,7,6,0,0,2),
    @(6,7,7,1,1,0),
    @(8,6,6,0,1,1),
    @(6,7,7,1,0,0),
    @(6,7,7,1,0,0),
    @(8,6,6,0,1,1),
    @(7,6,7,0,1,1),
    @(6,7,7,1,0,1),
    @(7,7,6,0,0,1),
    @(8,6,6,0,1,1),
    @(7,6,7,0,2,0),
    @(6,7,7,1,0,0),
    @(6,7,7,1,1,0),
    @(7,6,7,1,1,0),
    @(7,6,7,0,1,0),
    @(7,7,6,0,1,1),
    @(7,7,6,1,1,1),
    @(7,6,7,0,1,0)
)

$startRow = 4
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $vals = $data[$i]
    for ($c = 0; $c -lt $vals.Length; $c++) {
        $ws.Cells.Item($row, 2 + $c).Value = $vals[$c]
    }
}
